$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Threat Management")
$ws.Range("B2").Value = "TEST"
